$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 88: LeetCode 1967 - Delete Characters to Make Fancy String ---
# Seed formatting/styles by copying an existing 9-column (no Notes) data row.
$ws.Range("A85:I85").Copy()
$ws.Range("A88:I88").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(88, 1).Value = 1967
$ws.Cells.Item(88, 2).Value = "Delete Characters to Make Fancy String"
$ws.Cells.Item(88, 3).Value = "#string"
$ws.Cells.Item(88, 4).Value = "easy"
$ws.Cells.Item(88, 5).Value = 1
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 12
$ws.Cells.Item(88, 8).Value = 45859
$ws.Cells.Item(88, 9).Value = 45859
$ws.Rows.Item(88).RowHeight = 34

# --- Row 89: LeetCode 465 - Optimal Account Balancing ---
# Seed formatting/styles by copying an existing 10-column (with Notes) data row.
$ws.Range("A87:J87").Copy()
$ws.Range("A89:J89").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(89, 1).Value = 465
$ws.Cells.Item(89, 2).Value = "Optimal Account Balancing"
$ws.Cells.Item(89, 3).Value = "#array #backtracking #dynamic-programming "
$ws.Cells.Item(89, 4).Value = "hard"
$ws.Cells.Item(89, 5).Value = 0
$ws.Cells.Item(89, 6).Value = 1
$ws.Cells.Item(89, 7).Value = 90
$ws.Cells.Item(89, 8).Value = 45859
$ws.Cells.Item(89, 9).Value = 45859
$ws.Cells.Item(89, 10).Value = "？？dp 难！"
$ws.Rows.Item(89).RowHeight = 68

# --- Update the view state to match the committed workbook ---
$win = $wb.Windows.Item(1)
$win.ScrollRow = 86
$win.ScrollColumn = 2
try { $win.TopLeftCell = $ws.Range("B86") } catch {}
$ws.Range("H90").Select()
